$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = 2
$ws.Range("J25").Value = 2.03
$ws.Range("K25").Value = -1.37
$ws.Range("I26").Value = 2.12
$ws.Range("J26").Value = -1.28
$ws.Range("K26").Value = -1.28
$ws.Range("H27").Value = 2.09
$ws.Range("I27").Value = -1.31
$ws.Range("J27").Value = -1.31
$ws.Range("K27").Value = -0.01000000000000001
$ws.Range("G28").Value = 2.01
$ws.Range("H28").Value = -1.39
$ws.Range("I28").Value = -1.39
$ws.Range("J28").Value = -0.08999999999999997
$ws.Range("K28").Value = 0.11
$ws.Range("F29").Value = 1.95
$ws.Range("G29").Value = -1.45
$ws.Range("H29").Value = -1.45
$ws.Range("I29").Value = -0.15
$ws.Range("J29").Value = 0.05
$ws.Range("K29").Value = -0.15
$ws.Range("E30").Value = 1.93
$ws.Range("F30").Value = -1.47
$ws.Range("G30").Value = -1.47
$ws.Range("H30").Value = -0.17
$ws.Range("I30").Value = 0.03000000000000001
$ws.Range("J30").Value = -0.17
$ws.Range("K30").Value = -1.27
$ws.Range("D31").Value = 1.94
$ws.Range("E31").Value = -1.46
$ws.Range("F31").Value = -1.46
$ws.Range("G31").Value = -0.16
$ws.Range("H31").Value = 0.03999999999999999
$ws.Range("I31").Value = -0.16
$ws.Range("J31").Value = -1.26
$ws.Range("K31").Value = -0.16
$ws.Range("C32").Value = 1.95
$ws.Range("D32").Value = -1.45
$ws.Range("E32").Value = -1.45
$ws.Range("F32").Value = -0.15
$ws.Range("G32").Value = 0.05
$ws.Range("H32").Value = -0.15
$ws.Range("I32").Value = -1.25
$ws.Range("J32").Value = -0.15
$ws.Range("K32").Value = 0.05
$ws.Range("B33").Value = 1.96
$ws.Range("C33").Value = -1.44
$ws.Range("D33").Value = -1.44
$ws.Range("E33").Value = -0.14
$ws.Range("F33").Value = 0.05999999999999998
$ws.Range("G33").Value = -0.14
$ws.Range("H33").Value = -1.24
$ws.Range("I33").Value = -0.14
$ws.Range("J33").Value = 0.05999999999999998
$ws.Range("K33").Value = -0.34
$ws.Range("B34").Value = -1.67
$ws.Range("C34").Value = -1.67
$ws.Range("D34").Value = -0.37
$ws.Range("E34").Value = -0.17
$ws.Range("F34").Value = -0.37
$ws.Range("G34").Value = -1.47
$ws.Range("H34").Value = -0.37
$ws.Range("I34").Value = -0.17
$ws.Range("J34").Value = -0.5700000000000001
$ws.Range("K34").Value = -0.27
$ws.Range("B35").Value = -1.55
$ws.Range("C35").Value = -0.25
$ws.Range("D35").Value = -0.04999999999999999
$ws.Range("E35").Value = -0.25
$ws.Range("F35").Value = -1.35
$ws.Range("G35").Value = -0.25
$ws.Range("H35").Value = -0.04999999999999999
$ws.Range("I35").Value = -0.45
$ws.Range("J35").Value = -0.15
$ws.Range("K35").Value = -0.55
$ws.Range("B36").Value = -0.23
$ws.Range("C36").Value = -0.03
$ws.Range("D36").Value = -0.23
$ws.Range("E36").Value = -1.33
$ws.Range("F36").Value = -0.23
$ws.Range("G36").Value = -0.03
$ws.Range("H36").Value = -0.43
$ws.Range("I36").Value = -0.13
$ws.Range("J36").Value = -0.53
$ws.Range("K36").Value = 0.27
$ws.Range("B37").Value = 0.03999999999999999
$ws.Range("C37").Value = -0.16
$ws.Range("D37").Value = -1.26
$ws.Range("E37").Value = -0.16
$ws.Range("F37").Value = 0.03999999999999999
$ws.Range("G37").Value = -0.36
$ws.Range("H37").Value = -0.06000000000000002
$ws.Range("I37").Value = -0.46
$ws.Range("J37").Value = 0.34
$ws.Range("K37").Value = 0.03999999999999999
$ws.Range("B38").Value = -0.13
$ws.Range("C38").Value = -1.23
$ws.Range("D38").Value = -0.13
$ws.Range("E38").Value = 0.07000000000000002
$ws.Range("F38").Value = -0.33
$ws.Range("G38").Value = -0.02999999999999998
$ws.Range("H38").Value = -0.43
$ws.Range("I38").Value = 0.3700000000000001
$ws.Range("J38").Value = 0.07000000000000002
$ws.Range("K38").Value = -0.02999999999999998
$ws.Range("B39").Value = -1.24
$ws.Range("C39").Value = -0.14
$ws.Range("D39").Value = 0.06
$ws.Range("E39").Value = -0.34
$ws.Range("F39").Value = -0.04
$ws.Range("G39").Value = -0.44
$ws.Range("H39").Value = 0.36
$ws.Range("I39").Value = 0.06
$ws.Range("J39").Value = -0.04
$ws.Range("K39").Value = -0.44
$ws.Range("B40").Value = 0.07000000000000001
$ws.Range("C40").Value = 0.27
$ws.Range("D40").Value = -0.13
$ws.Range("E40").Value = 0.17
$ws.Range("F40").Value = -0.23
$ws.Range("G40").Value = 0.5700000000000001
$ws.Range("H40").Value = 0.27
$ws.Range("I40").Value = 0.17
$ws.Range("J40").Value = -0.23
$ws.Range("K40").Value = -0.63
$ws.Range("B41").Value = 0.14
$ws.Range("C41").Value = -0.26
$ws.Range("D41").Value = 0.04000000000000004
$ws.Range("E41").Value = -0.36
$ws.Range("F41").Value = 0.4400000000000001
$ws.Range("G41").Value = 0.14
$ws.Range("H41").Value = 0.04000000000000004
$ws.Range("I41").Value = -0.36
$ws.Range("J41").Value = -0.76
$ws.Range("K41").Value = -0.76
$ws.Range("B42").Value = -0.26
$ws.Range("C42").Value = 0.04000000000000001
$ws.Range("D42").Value = -0.36
$ws.Range("E42").Value = 0.4400000000000001
$ws.Range("F42").Value = 0.14
$ws.Range("G42").Value = 0.04000000000000001
$ws.Range("H42").Value = -0.36
$ws.Range("I42").Value = -0.76
$ws.Range("J42").Value = -0.76
$ws.Range("K42").Value = 0.9400000000000001
$ws.Range("B43").Value = 0.1
$ws.Range("C43").Value = -0.3
$ws.Range("D43").Value = 0.5
$ws.Range("E43").Value = 0.2
$ws.Range("F43").Value = 0.1
$ws.Range("G43").Value = -0.3
$ws.Range("H43").Value = -0.7000000000000001
$ws.Range("I43").Value = -0.7000000000000001
$ws.Range("J43").Value = 1
$ws.Range("K43").Value = -0.4
$ws.Range("B44").Value = -0.09000000000000002
$ws.Range("C44").Value = 0.71
$ws.Range("D44").Value = 0.41
$ws.Range("E44").Value = 0.31
$ws.Range("F44").Value = -0.09000000000000002
$ws.Range("G44").Value = -0.49
$ws.Range("H44").Value = -0.49
$ws.Range("I44").Value = 1.21
$ws.Range("J44").Value = -0.19
$ws.Range("B45").Value = 0.6200000000000001
$ws.Range("C45").Value = 0.3200000000000001
$ws.Range("D45").Value = 0.2200000000000001
$ws.Range("E45").Value = -0.1799999999999999
$ws.Range("F45").Value = -0.58
$ws.Range("G45").Value = -0.58
$ws.Range("H45").Value = 1.12
$ws.Range("I45").Value = -0.2799999999999999
$ws.Range("B46").Value = 0.11
$ws.Range("C46").Value = 0.009999999999999986
$ws.Range("D46").Value = -0.39
$ws.Range("E46").Value = -0.79
$ws.Range("F46").Value = -0.79
$ws.Range("G46").Value = 0.91
$ws.Range("H46").Value = -0.49
$ws.Range("B47").Value = 0.06000000000000001
$ws.Range("C47").Value = -0.34
$ws.Range("D47").Value = -0.74
$ws.Range("E47").Value = -0.74
$ws.Range("F47").Value = 0.9600000000000001
$ws.Range("G47").Value = -0.44
$ws.Range("B48").Value = -0.4
$ws.Range("C48").Value = -0.8
$ws.Range("D48").Value = -0.8
$ws.Range("E48").Value = 0.9
$ws.Range("F48").Value = -0.5
$ws.Range("B49").Value = -0.79
$ws.Range("C49").Value = -0.79
$ws.Range("D49").Value = 0.91
$ws.Range("E49").Value = -0.49
$ws.Range("B50").Value = -0.67
$ws.Range("C50").Value = 1.03
$ws.Range("D50").Value = -0.37
$ws.Range("B51").Value = 1
$ws.Range("C51").Value = -0.4
$ws.Range("B52").Value = -0.53
